$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.619.05'
$ws.Range("E2").Value = '  -0.39%  '

$ws.Range("D3").Value = '2.667.32'
$ws.Range("E3").Value = '  -1.09%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''598.71'
$ws.Range("E5").Value = '  -1.69%  '

$ws.Range("D6").Value = '''156.46'
$ws.Range("E6").Value = '  -0.97%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '''0.616'
$ws.Range("E8").Value = '  +4.92%  '

$ws.Range("D9").Value = '''0.129'
$ws.Range("E9").Value = '  +4.11%  '

$ws.Range("D10").Value = '''0.399'
$ws.Range("E10").Value = '  -0.97%  '

$ws.Range("D11").Value = '''5.85'
$ws.Range("E11").Value = '  -2.86%  '

$ws.Range("E12").Value = '  -0.24%  '

$ws.Range("D13").Value = '''29.24'
$ws.Range("E13").Value = '  -3.70%  '

$ws.Range("D14").Value = '''0.0000196'
$ws.Range("E14").Value = '  -2.62%  '

$ws.Range("D15").Value = '3.146.32'
$ws.Range("E15").Value = '  -1.20%  '

$ws.Range("D16").Value = '65.441.16'
$ws.Range("E16").Value = '  -0.43%  '

$ws.Range("D17").Value = '2.663.40'
$ws.Range("E17").Value = '  -1.70%  '

$ws.Range("D18").Value = '''12.79'
$ws.Range("E18").Value = '  +1.11%  '

$ws.Range("D19").Value = '''4.77'
$ws.Range("E19").Value = '  -2.45%  '

$ws.Range("D20").Value = '''7.54'
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").Value = '''351.06'
$ws.Range("E21").Value = '  -2.36%  '

$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").Value = '''69.62'
$ws.Range("E23").Value = '  -1.55%  '

$ws.Range("D24").Value = '''0.0000111'
$ws.Range("E24").Value = '  +3.60%  '

$ws.Range("D25").Value = '''9.58'
$ws.Range("E25").Value = '  -2.71%  '

$ws.Range("D26").Value = '''1.63'
$ws.Range("E26").Value = '  -1.84%  '

$ws.Range("D27").Value = '''0.167'
$ws.Range("E27").Value = '  -2.73%  '

$ws.Range("D28").Value = '''1.58'
$ws.Range("E28").Value = '  -5.73%  '

$ws.Range("D29").Value = '''7.99'
$ws.Range("E29").Value = '  -4.71%  '

$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("D31").Value = '''2.13'
$ws.Range("E31").Value = '  -3.68%  '

$ws.Range("D32").Value = '''522.59'
$ws.Range("E32").Value = '  -3.40%  '

$ws.Range("D33").Value = '''1.76'
$ws.Range("E33").Value = '  -2.27%  '

$ws.Range("D34").Value = '''6.43'
$ws.Range("E34").Value = '  -3.80%  '

$ws.Range("D35").Value = '''5.44'
$ws.Range("E35").Value = '  +0.89%  '

$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").Value = '''0.421'
$ws.Range("E36").Value = '  -2.57%  '

$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").Value = '''20.54'
$ws.Range("E37").Value = '  -1.51%  '

$ws.Range("D38").Value = '''0.999'
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("D39").Value = '''157.62'
$ws.Range("E39").Value = '  -3.53%  '

$ws.Range("D40").Value = '''1.93'
$ws.Range("E40").Value = '  -3.42%  '

$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("D42").Value = '''163.03'
$ws.Range("E42").Value = '  -5.49%  '

$ws.Range("D43").Value = '''4.11'
$ws.Range("E43").Value = '  -1.94%  '

$ws.Range("D44").Value = '''2.29'
$ws.Range("E44").Value = '  +1.09%  '

$ws.Range("D45").Value = '''0.0606'
$ws.Range("E45").Value = '  -1.49%  '

$ws.Range("D46").Value = '''22.67'
$ws.Range("E46").Value = '  -3.93%  '

$ws.Range("E47").Value = '  -2.99%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0263'
$ws.Range("E48").Value = '  +13.98%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '''0.0256'
$ws.Range("E49").Value = '  -3.78%  '

$ws.Range("D50").Value = '''0.0997'
$ws.Range("E50").Value = '  +0.40%  '

$ws.Range("D51").Value = '''20.04'
$ws.Range("E51").Value = '  -4.89%  '
